$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PE, SEPTEMBER")

# Row 12
$ws.Range("C12").Value = 45924
$ws.Range("G12").Value = 517936275
$ws.Range("I12").Formula = "=1330452-56255.04"

# Row 13
$ws.Range("C13").Value = 45924
$ws.Range("G13").Value = 517936283
$ws.Range("I13").Formula = "=1283160-53990.4"

# Row 14
$ws.Range("C14").Value = 45924
$ws.Range("G14").Formula = "=517935815"
$ws.Range("I14").Formula = "=1324098-53803.62"

# Row 15
$ws.Range("C15").Value = 45924
$ws.Range("G15").Value = 517935830
$ws.Range("I15").Formula = "=1366356-56604.96"

# Update the active cell selection to match the saved file state
$ws.Range("C15").Select()

$wb.Save()
